$d = $word.ActiveDocument

$lastPara = $d.Paragraphs.Last
$endRange = $lastPara.Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$insertionPara = $d.Paragraphs.Last
$insertionRange = $insertionPara.Range

$xml = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
<w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:bidi w:val="0"/>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:b/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:bidi w:val="0"/>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:b/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>2.2 Initial feasibility study</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:bidi w:val="0"/>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:b/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:bidi w:val="0"/>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:b/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>2.2.1 Lack of capital</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:bidi w:val="0"/>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:b/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:bidi w:val="0"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:b/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b w:val="false"/>
          <w:bCs w:val="false"/>
        </w:rPr>
        <w:t>C</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b w:val="false"/>
          <w:bCs w:val="false"/>
        </w:rPr>
        <w:t>apital is the money that a bank has obtained from its shareholders and other investors and any profit that it has made and not paid out</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:bidi w:val="0"/>
        <w:jc w:val="both"/>
        <w:rPr/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b w:val="false"/>
          <w:bCs w:val="false"/>
        </w:rPr>
        <w:t xml:space="preserve">Capital is a key ingredient for safe and sound banking. </w:t>
      </w:r>
      <w:r>
        <w:rPr/>
        <w:t>Banks take on risks and may suffer losses if the risks. To stay safe and protect people’s deposits, banks have to be able to absorb such losses and keep going in good times and bad. That’s what bank capital is used for.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:bidi w:val="0"/>
        <w:jc w:val="both"/>
        <w:rPr/>
      </w:pPr>
      <w:r>
        <w:rPr/>
        <w:t>“</w:t>
      </w:r>
      <w:r>
        <w:rPr/>
        <w:t xml:space="preserve">Rajshahi Krishi Unnayan Bank” is a government organization. So, government fixes the budget of the capital for the bank. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:bidi w:val="0"/>
        <w:jc w:val="both"/>
        <w:rPr/>
      </w:pPr>
      <w:r>
        <w:rPr/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:bidi w:val="0"/>
        <w:jc w:val="both"/>
        <w:rPr/>
      </w:pPr>
      <w:r>
        <w:rPr/>
        <w:t xml:space="preserve">The main function of the bank is to give loan to the consumers, </w:t>
      </w:r>
      <w:r>
        <w:rPr/>
        <w:t xml:space="preserve">handling deposits of the consumers. The bank has a vast number of consumers. The authorized capital of the bank is 10000 million dollar. The bank was paid up 825 million dollar. Sometimes the </w:t>
      </w:r>
      <w:r>
        <w:rPr/>
        <w:t xml:space="preserve">transactions is very large. The capital can not support the transactions. As a result, bank has to take loans from or make a financial request to the Bangladesh Bank. This causes the system to suffer losses. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:bidi w:val="0"/>
        <w:jc w:val="left"/>
        <w:rPr/>
      </w:pPr>
      <w:r>
        <w:rPr/>
      </w:r>
    </w:p>
    
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$insertionRange.InsertXML($xml)

$countAfterInsert = $d.Paragraphs.Count
$lastRealPara = $d.Paragraphs.Item($countAfterInsert - 1)
$cleanupRange = $d.Range($lastRealPara.Range.End - 1, $d.Paragraphs.Last.Range.End)
$cleanupRange.Delete()

Write-Output $d.Paragraphs.Count
